# Apply crypto price/volume updates to Sheet1 (cryptos.xlsx)
# Leading apostrophe forces text interpretation so numeric-looking
# strings (e.g. "1.340", "0.9983") are not auto-converted to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.072.45"
$ws.Range("E2").Value = "'  -0.14%  "

$ws.Range("D3").Value = "'1.879.90"

$ws.Range("D4").Value = "'0.9983"
$ws.Range("E4").Value = "'  -0.09%  "

$ws.Range("D5").Value = "'243.20"
$ws.Range("E5").Value = "'  -3.69%  "

$ws.Range("D6").Value = "'0.9980"
$ws.Range("E6").Value = "'  -0.14%  "

$ws.Range("D7").Value = "'0.4917"
$ws.Range("E7").Value = "'  -3.29%  "

$ws.Range("D8").Value = "'0.2944"
$ws.Range("E8").Value = "'  -2.25%  "

$ws.Range("D9").Value = "'0.06610"
$ws.Range("E9").Value = "'  -3.17%  "

$ws.Range("D10").Value = "'1.880.28"
$ws.Range("E10").Value = "'  -1.33%  "

$ws.Range("D11").Value = "'16.67"
$ws.Range("E11").Value = "'  -3.70%  "

$ws.Range("D12").Value = "'0.07179"
$ws.Range("E12").Value = "'  -1.96%  "

$ws.Range("D13").Value = "'0.6662"
$ws.Range("E13").Value = "'  -3.87%  "

$ws.Range("D14").Value = "'86.33"
$ws.Range("E14").Value = "'  -0.90%  "

$ws.Range("D15").Value = "'4.889"
$ws.Range("E15").Value = "'  -0.46%  "

$ws.Range("D16").Value = "'30.025.88"
$ws.Range("E16").Value = "'  -0.30%  "

$ws.Range("D17").Value = "'0.000007803"
$ws.Range("E17").Value = "'  -7.04%  "

$ws.Range("D18").Value = "'0.9977"
$ws.Range("E18").Value = "'  -0.16%  "

$ws.Range("D19").Value = "'12.79"
$ws.Range("E19").Value = "'  -2.19%  "

$ws.Range("D20").Value = "'2.122.45"
$ws.Range("E20").Value = "'  -1.30%  "

$ws.Range("D21").Value = "'0.9991"
$ws.Range("E21").Value = "'  +0.08%  "

$ws.Range("D22").Value = "'4.775"
$ws.Range("E22").Value = "'  -0.94%  "

$ws.Range("D23").Value = "'5.846"
$ws.Range("E23").Value = "'  +1.69%  "

$ws.Range("D24").Value = "'9.085"
$ws.Range("E24").Value = "'  -2.56%  "

$ws.Range("D25").Value = "'151.44"
$ws.Range("E25").Value = "'  +2.59%  "

$ws.Range("D26").Value = "'142.77"
$ws.Range("E26").Value = "'  +6.06%  "

$ws.Range("D27").Value = "'16.96"
$ws.Range("E27").Value = "'  -1.25%  "

$ws.Range("D28").Value = "'1.894"
$ws.Range("E28").Value = "'  -5.50%  "

$ws.Range("D29").Value = "'1.385"
$ws.Range("E29").Value = "'  -1.06%  "

$ws.Range("D30").Value = "'4.201"
$ws.Range("E30").Value = "'  -2.32%  "

$ws.Range("D31").Value = "'0.08772"
$ws.Range("E31").Value = "'  -1.07%  "

$ws.Range("D32").Value = "'3.981"
$ws.Range("E32").Value = "'  -0.62%  "

$ws.Range("D33").Value = "'0.05014"
$ws.Range("E33").Value = "'  -1.10%  "

$ws.Range("D34").Value = "'0.7244"
$ws.Range("E34").Value = "'  +0.20%  "

$ws.Range("D35").Value = "'1.113"
$ws.Range("E35").Value = "'  -2.89%  "

$ws.Range("D36").Value = "'2.662"
$ws.Range("E36").Value = "'  -1.11%  "

$ws.Range("D37").Value = "'0.01797"
$ws.Range("E37").Value = "'  +6.06%  "

$ws.Range("D38").Value = "'2.687"
$ws.Range("E38").Value = "'  -4.64%  "

$ws.Range("D39").Value = "'2.157"
$ws.Range("E39").Value = "'  -6.02%  "

$ws.Range("D40").Value = "'0.9377"
$ws.Range("E40").Value = "'  -2.89%  "

$ws.Range("D41").Value = "'0.9980"
$ws.Range("E41").Value = "'  -0.11%  "

$ws.Range("B42").Value = "'FraxShare"
$ws.Range("C42").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'5.752"
$ws.Range("E42").Value = "'  -5.50%  "

$ws.Range("B43").Value = "'TheSandbox"
$ws.Range("C43").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4229"
$ws.Range("E43").Value = "'  -1.86%  "

$ws.Range("B44").Value = "'Quant"
$ws.Range("C44").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'103.34"
$ws.Range("E44").Value = "'  -1.68%  "

$ws.Range("D45").Value = "'7.355"
$ws.Range("E45").Value = "'  -4.36%  "

$ws.Range("D46").Value = "'0.1270"
$ws.Range("E46").Value = "'  -0.97%  "

$ws.Range("E47").Value = "'  -1.17%  "

$ws.Range("D48").Value = "'32.74"
$ws.Range("E48").Value = "'  -2.05%  "

$ws.Range("D49").Value = "'8.296"
$ws.Range("E49").Value = "'  -1.40%  "

$ws.Range("D50").Value = "'0.3769"
$ws.Range("E50").Value = "'  -1.29%  "

$ws.Range("B51").Value = "'NEARProtocol"
$ws.Range("C51").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.340"
$ws.Range("E51").Value = "'  -2.06%  "
